# Apply the weekly Fruta/hortaliza data shuffle to rows 2-9 (columns D, M, Q, S, T)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) - swap/shuffle serial date values across rows 2-9
$ws.Range("D2").Value = 44330
$ws.Range("D3").Value = 44316
$ws.Range("D4").Value = 44309
$ws.Range("D5").Value = 44306
$ws.Range("D6").Value = 44313
$ws.Range("D7").Value = 44327
$ws.Range("D8").Value = 44302
$ws.Range("D9").Value = 44322

# Column M (Volumen) - shuffle values across rows 2-9
$ws.Range("M2").Value = 60
$ws.Range("M3").Value = 120
$ws.Range("M4").Value = 80
$ws.Range("M6").Value = 120
$ws.Range("M8").Value = 80
$ws.Range("M9").Value = 60

# Rows 4 and 5: swap unidad de comercialización / precio $/kg / kg-unidad
$ws.Range("Q4").Value = "$/caja 14 kilos granel"
$ws.Range("S4").Value = 821
$ws.Range("T4").Value = 14

$ws.Range("Q5").Value = "$/caja 10 kilos empedrada"
$ws.Range("S5").Value = 11500
$ws.Range("T5").Value = 1
